$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting all existing data right by one.
$ws.Range("A1").EntireColumn.Insert()

# New column width for column A (narrow column for the lab number code).
$ws.Columns.Item(1).ColumnWidth = 6.86

# New header + laboratory number value.
$ws.Range("A1").Value = "Lab. #"
$ws.Range("A2").Value = 7184

# Highlight the whole data row (A2:U2) with the light-green fill used for the
# added laboratory-number row, matching the new cellXfs style (fillId=2).
$dataRow = $ws.Range("A2:U2")
$dataRow.Interior.Color = 12379352
